$d = $word.ActiveDocument

# 1) The "_GoBack" bookmark (Word's "last edit location" marker) needs to
#    move from its old spot - the indent-only paragraph right after
#    "Rétablir le robot s'il tombe ?" - to the new paragraph we're about to
#    add below. Remove the stale one first so we don't end up with two
#    bookmarks sharing the name "_GoBack" once the new one is inserted.
#    (It's a hidden bookmark, so it won't show in Bookmarks.Count but is
#    still reachable by name.)
try {
    $old = $d.Bookmarks.Item("_GoBack")
    $old.Delete()
} catch {
}

# 2) Fill in the empty paragraph right after "Autres joyeusetés" with the
#    italic quoted note ("Etude des problématiques ..."), reproducing the
#    exact run/proofErr/bookmark layout from the authored OOXML.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "`r" -or $p.Range.Text -eq "") {
        $prev = $p.Previous()
        if ($prev -ne $null -and $prev.Range.Text.Trim() -eq "Autres joyeusetés") {
            $target = $p
            break
        }
    }
}

if ($target -eq $null) {
    throw "Could not locate the empty placeholder paragraph after 'Autres joyeusetés'"
}

$xml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:i/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:i/>
    </w:rPr>
    <w:t>« </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:i/>
    </w:rPr>
    <w:t xml:space="preserve">Etude des </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:i/>
    </w:rPr>
    <w:t>problématiques</w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
  <w:r>
    <w:rPr>
      <w:i/>
    </w:rPr>
    <w:t xml:space="preserve"> associées à un robot marcheur à 4 pattes : maintenir l'équilibre, porter une masse, se relever en cas de chute, etc. Réaliser un court état de l'art des robots marcheurs, et réfléchir à ce qu'on pourrait faire avec le robot disponible au Dpt, le </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:i/>
    </w:rPr>
    <w:t>Quattro</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:i/>
    </w:rPr>
    <w:t xml:space="preserve"> de Roboticia.</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:i/>
    </w:rPr>
    <w:t> »</w:t>
  </w:r>
</w:p>
"@

$target.Range.InsertXML($xml) | Out-Null
